$d = $word.ActiveDocument

$replacements = @(
    @("809×7=", "304×3="),
    @("516×2=", "434×5="),
    @("174×7=", "343×4="),
    @("419×3=", "358×9="),
    @("799×9=", "849×8="),
    @("548×7=", "591×9="),
    @("188×7=", "899×6="),
    @("352×3=", "207×8="),
    @("543×5=", "333×8="),
    @("222×6=", "744×5="),
    @("422×6=", "862×8="),
    @("736×5=", "558×2="),
    @("691×9=", "492×9="),
    @("439×8=", "243×9="),
    @("161×9=", "940×3="),
    @("850×2=", "440×3="),
    @("297×9=", "516×7="),
    @("929×9=", "236×5="),
    @("147×3=", "426×8="),
    @("914×6=", "841×8="),
    @("596×8=", "745×3="),
    @("139×5=", "330×6="),
    @("750×3=", "791×6="),
    @("180×7=", "650×4="),
    @("520×6=", "502×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
